$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column before C ("Validación") -------------------------
$ws.Columns("C").Insert()

# Give the new column the same width as A/B, and make sure it (and its
# header/data cells) use a text number format, same shape as the other
# "text" columns in this report (preserves things like leading zeros).
$ws.Columns("C").ColumnWidth = 28.19
$ws.Range("C1:C6").NumberFormat = "@"

# Header for the new column
$ws.Range("C1").Value = "Validación"

# Rows 2-6 in column C stay blank (no value) - nothing else to do there.

# --- Fix up "No. Recibo sin PCC" (A) / "Key" (B) on rows 2 and 4 ---------
# These two rows had a typo'd / truncated receipt number; correct them to
# match the rest of the rows for the same sale. Force text first so the
# long numeric string round-trips exactly (no scientific notation /
# precision loss), then copy the original cell formatting back on top so
# the cell keeps its normal (non-text-numfmt) style, same as the diff.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "50000000210344998"
$ws.Range("B2").Value = "50000000210344998EDS3234"
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A2").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "50000000210344998"
$ws.Range("B4").Value = "50000000210344998EDS3234"
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A4").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Reset the view (no frozen/scrolled selection) -----------------------
$ws.Range("A1").Select()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
